# Update the "想去人数" (interested-people count) figures on both the
# "展览" and "全部类型" sheets to reflect the refreshed scrape data.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 672
    $ws.Range("F3").Value = 3949
    $ws.Range("F5").Value = 734
}
